$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.875.06"
$ws.Range("D3").Value = "1.752.39"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'236.23"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.5147"
$ws.Range("E7").Value = "  +5.13%  "
$ws.Range("D8").Value = "'40.46"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").Value = "'0.2678"
$ws.Range("E9").Value = "  +5.25%  "
$ws.Range("D10").Value = "'0.06180"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "1.773.75"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").Value = "'0.06945"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "'15.46"
$ws.Range("E13").Value = "  +4.87%  "
$ws.Range("D14").Value = "'0.6372"
$ws.Range("E14").Value = "  +14.15%  "
$ws.Range("D15").Value = "'4.493"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "'77.99"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "'0.9981"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'0.9982"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "25.889.81"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "'11.63"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").Value = "'0.000006669"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "1.987.15"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").Value = "'4.064"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "'8.285"
$ws.Range("E24").Value = "  +4.97%  "
$ws.Range("D25").Value = "'5.175"
$ws.Range("E25").Value = "  +3.69%  "
$ws.Range("D26").Value = "'136.17"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'15.09"
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("D29").Value = "'1.762"
$ws.Range("E29").Value = "  -3.19%  "
$ws.Range("D30").Value = "'102.73"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("D32").Value = "'3.684"
$ws.Range("D33").Value = "'3.392"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").Value = "'2.636"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "'0.9997"
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("E37").Value = "  +2.92%  "
$ws.Range("D38").Value = "'2.738"
$ws.Range("E38").Value = "  +3.88%  "
$ws.Range("D39").Value = "'0.01564"
$ws.Range("E39").Value = "  +4.44%  "
$ws.Range("D40").Value = "'1.938"
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("D41").Value = "'0.9988"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'102.13"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").Value = "'0.3841"
$ws.Range("E43").Value = "  +3.57%  "
$ws.Range("D44").Value = "'0.7490"
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("D45").Value = "'4.903"
$ws.Range("E45").Value = "  -4.33%  "
$ws.Range("E46").Value = "  +5.40%  "
$ws.Range("D47").Value = "'0.1104"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").Value = "'5.983"
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("D49").Value = "'30.11"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "'52.61"
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("E51").Value = "  +0.42%  "
